$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: shift values so C1=prediction, D1=rejection-f, E1=max
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Row 2: C becomes the species string, E becomes a numeric value
$ws.Range("C2").Value = "s__CAG-631 sp000433015"
$ws.Range("E2").Value = 0.981899820276805

# Row 3: C becomes the species string, E becomes a numeric value
$ws.Range("C3").Value = "s__CAG-631 sp000433015"
$ws.Range("E3").Value = 0.9811743162699852
